$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both carry the same "想去人数" (F column)
# counts and both need to be bumped identically.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 131
    $ws.Range("F6").Value = 458
    $ws.Range("F7").Value = 151
    $ws.Range("F8").Value = 72
    $ws.Range("F9").Value = 586
}
